# Auto-generated from the OOXML unified diff for Tiamat_Profits workbook.
# Updates cached market-price / profit values across the ALC, ARM, BSM, CUL,
# GSM and LTW sheets; also removes the stale H:N figures on GSM rows 125-141
# (the leves there moved down one row) and adds the corresponding H:N figures
# on LTW rows 124-141.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# ALC: set updated values
$ws.Range("H129").Value = 1893.8695
$ws.Range("J129").Value = 2129.0264
$ws.Range("L129").Value = 6387.0792
$ws.Range("N129").Value = -16387.0792
$ws.Range("H137").Value = 1406.5625
$ws.Range("I137").Value = 1207.0435
$ws.Range("J137").Value = 1590.12
$ws.Range("K137").Value = 3621.1305
$ws.Range("L137").Value = 4770.36
$ws.Range("M137").Value = -1071.1305
$ws.Range("N137").Value = -9870.360000000001
$ws.Range("H141").Value = 2545.25
$ws.Range("I141").Value = 2034.6428
$ws.Range("J141").Value = 3736.6667
$ws.Range("K141").Value = 6103.928400000001
$ws.Range("L141").Value = 11210.0001
$ws.Range("M141").Value = -923.9284000000007
$ws.Range("N141").Value = -21570.0001

$ws = $wb.Worksheets.Item("ARM")

# ARM: set updated values
$ws.Range("H32").Value = 172518.55
$ws.Range("I32").Value = 184529.97
$ws.Range("J32").Value = 105909.73
$ws.Range("K32").Value = 184529.97
$ws.Range("L32").Value = 105909.73
$ws.Range("M32").Value = -184242.97
$ws.Range("N32").Value = -106483.73
$ws.Range("H45").Value = 1753.2333
$ws.Range("I45").Value = 1534.0952
$ws.Range("K45").Value = 1534.0952
$ws.Range("M45").Value = -1157.0952
$ws.Range("H74").Value = 59142.516
$ws.Range("I74").Value = 66339.61
$ws.Range("J74").Value = 3365
$ws.Range("K74").Value = 66339.61
$ws.Range("L74").Value = 3365
$ws.Range("M74").Value = -65465.61
$ws.Range("N74").Value = -5113
$ws.Range("H77").Value = 59142.516
$ws.Range("I77").Value = 66339.61
$ws.Range("J77").Value = 3365
$ws.Range("K77").Value = 331698.05
$ws.Range("L77").Value = 16825
$ws.Range("M77").Value = -327330.05
$ws.Range("N77").Value = -25561
$ws.Range("H132").Value = 2493469.8
$ws.Range("I132").Value = 3644223.2
$ws.Range("J132").Value = 959131.5
$ws.Range("K132").Value = 10932669.6
$ws.Range("L132").Value = 2877394.5
$ws.Range("M132").Value = -10930139.6
$ws.Range("N132").Value = -2882454.5

$ws = $wb.Worksheets.Item("BSM")

# BSM: set updated values
$ws.Range("H64").Value = 270.14285
$ws.Range("I64").Value = 184
$ws.Range("J64").Value = 304.6
$ws.Range("K64").Value = 184
$ws.Range("L64").Value = 304.6
$ws.Range("M64").Value = 41
$ws.Range("N64").Value = -754.6
$ws.Range("H67").Value = 270.14285
$ws.Range("I67").Value = 184
$ws.Range("J67").Value = 304.6
$ws.Range("K67").Value = 184
$ws.Range("L67").Value = 304.6
$ws.Range("M67").Value = 596
$ws.Range("N67").Value = -1864.6
$ws.Range("H105").Value = 1779.091
$ws.Range("I105").Value = 1567.1428
$ws.Range("J105").Value = 2150
$ws.Range("K105").Value = 1567.1428
$ws.Range("L105").Value = 2150
$ws.Range("M105").Value = 179.8571999999999
$ws.Range("N105").Value = -5644

$ws = $wb.Worksheets.Item("CUL")

# CUL: set updated values
$ws.Range("H68").Value = 3087.4822
$ws.Range("I68").Value = 10720.2
$ws.Range("J68").Value = 1428.1957
$ws.Range("K68").Value = 32160.6
$ws.Range("L68").Value = 4284.5871
$ws.Range("M68").Value = -31349.6
$ws.Range("N68").Value = -5906.5871
$ws.Range("H71").Value = 3087.4822
$ws.Range("I71").Value = 10720.2
$ws.Range("J71").Value = 1428.1957
$ws.Range("K71").Value = 96481.8
$ws.Range("L71").Value = 12853.7613
$ws.Range("M71").Value = -92425.8
$ws.Range("N71").Value = -20965.7613
$ws.Range("H87").Value = 149144.28
$ws.Range("I87").Value = 3498.5
$ws.Range("J87").Value = 343338.66
$ws.Range("K87").Value = 10495.5
$ws.Range("L87").Value = 1030015.98
$ws.Range("M87").Value = -9247.5
$ws.Range("N87").Value = -1032511.98
$ws.Range("H90").Value = 149144.28
$ws.Range("I90").Value = 3498.5
$ws.Range("J90").Value = 343338.66
$ws.Range("K90").Value = 31486.5
$ws.Range("L90").Value = 3090047.94
$ws.Range("M90").Value = -25246.5
$ws.Range("N90").Value = -3102527.94

$ws = $wb.Worksheets.Item("GSM")

# GSM: set updated values
$ws.Range("H70").Value = 4036.2856
$ws.Range("I70").Value = 3953.4736
$ws.Range("K70").Value = 3953.4736
$ws.Range("M70").Value = -3683.4736
$ws.Range("H73").Value = 4036.2856
$ws.Range("I73").Value = 3953.4736
$ws.Range("K73").Value = 3953.4736
$ws.Range("M73").Value = -3017.4736

# GSM: clear cells removed in the target revision
$clearCells = @("H125", "I125", "J125", "K125", "L125", "N125", "H126", "I126", "J126", "K126", "L126", "M126", "N126", "H127", "I127", "J127", "K127", "L127", "N127", "H128", "I128", "J128", "K128", "L128", "N128", "H129", "I129", "J129", "K129", "L129", "N129", "H130", "I130", "J130", "K130", "L130", "N130", "H131", "I131", "J131", "K131", "L131", "M131", "N131", "H132", "I132", "J132", "K132", "L132", "M132", "N132", "H133", "I133", "J133", "K133", "L133", "N133", "H134", "I134", "J134", "K134", "L134", "H135", "I135", "J135", "K135", "L135", "N135", "H136", "I136", "J136", "K136", "L136", "N136", "H137", "I137", "J137", "K137", "L137", "N137", "H138", "I138", "J138", "K138", "L138", "H139", "I139", "J139", "K139", "L139", "H140", "I140", "J140", "K140", "L140", "N140", "H141", "I141", "J141", "K141", "L141", "N141")
foreach ($c in $clearCells) {
    $ws.Range($c).ClearContents()
}

$ws = $wb.Worksheets.Item("LTW")

# LTW: set updated values
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H125").Value = 40000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 40000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 39475
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 39475
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 39475
$ws.Range("N130").Value = -49515
$ws.Range("H131").Value = 30560
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 30560
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 30560
$ws.Range("N131").Value = -40640
$ws.Range("H132").Value = 318119.22
$ws.Range("I132").Value = 96642.48
$ws.Range("J132").Value = 591708.1
$ws.Range("K132").Value = 289927.44
$ws.Range("L132").Value = 1775124.3
$ws.Range("M132").Value = -287397.44
$ws.Range("N132").Value = -1780184.3
$ws.Range("H133").Value = 29663
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 29663
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 29663
$ws.Range("N133").Value = -34723
$ws.Range("H134").Value = 40000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 40000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 40000
$ws.Range("N134").Value = -50140
$ws.Range("H135").Value = 49972.727
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 49972.727
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 49972.727
$ws.Range("N135").Value = -60112.727
$ws.Range("H136").Value = 323814.25
$ws.Range("I136").Value = 455262.88
$ws.Range("J136").Value = 2495.4443
$ws.Range("K136").Value = 1365788.64
$ws.Range("L136").Value = 7486.3329
$ws.Range("M136").Value = -1363238.64
$ws.Range("N136").Value = -12586.3329
$ws.Range("H137").Value = 39879.832
$ws.Range("I137").Value = 40000
$ws.Range("J137").Value = 39874.61
$ws.Range("K137").Value = 40000
$ws.Range("L137").Value = 39874.61
$ws.Range("M137").Value = -34900
$ws.Range("N137").Value = -50074.61
$ws.Range("H138").Value = 41986.668
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 41986.668
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 41986.668
$ws.Range("N138").Value = -52266.668
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 51250
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 51250
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 51250
$ws.Range("N140").Value = -61610
$ws.Range("H141").Value = 69715
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 69715
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 69715
$ws.Range("N141").Value = -80075

